$wb = $excel.ActiveWorkbook

# Worksheets.Item(<name>) resolves case-insensitively, and this workbook has
# both "Vector_bf" and "Vector_BF" sheets, so address every sheet by its
# (unambiguous) 1-based position instead of by name.
$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsVecbf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)   # Vector_BF
$wsVecAlpha = $wb.Worksheets.Item(7)   # Vector_Alpha

# All of the cells below already hold plain-text (shared-string) values even
# though some of them look like numbers (e.g. "4.2", "0", "-1 + ...y").
# Assigning a numeric-looking string straight to .Value lets COM infer a
# number type and silently re-type the cell, so force the cell to Text
# first, write the literal string, then restore the default ("Normal")
# style so no stray formatting is left behind.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $wsFollower.Range("A2") "-1 + 2.3126160715899395y"
Set-TextValue $wsFollower.Range("B2") "6.631633036246801"
Set-TextValue $wsFollower.Range("E2") "2.3000000000000003"
Set-TextValue $wsFollower.Range("F2") "1.5"
Set-TextValue $wsFollower.Range("A3") "-1 + 0.35138390741670333y"
Set-TextValue $wsFollower.Range("B3") "0.15956689447512118"
Set-TextValue $wsFollower.Range("E3") "0"
Set-TextValue $wsFollower.Range("F3") "7.4"

Set-TextValue $wsVecbf.Range("A2") "-44.5597043627674"

Set-TextValue $wsVecBF.Range("A2") "2.2008662110497585"
Set-TextValue $wsVecBF.Range("A3") "-112.43635566709146"

# Vector_Alpha!A2 is a genuine numeric cell (no shared-string / text marker
# in the source), so assign it as a real number.
$wsVecAlpha.Range("A2").Value = 0.6486160925832967
